$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.83761066198349
$ws.Range("B1").Value = 1.083747863769531
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.137610673904419
$ws.Range("E1").Value = 1.081295251846313
